$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF column holds the game-date string for each row (BF2:BF31). The
# original text "5-22-2012-13" was off by a day because of how the NBA
# stats site displayed the date; replace it with the corrected
# "2013-05-22" value. Force the cell to Text first so Excel doesn't
# auto-convert the literal into a date serial number, then restore the
# "Normal" cell style so no stray per-cell number format sticks around.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2013-05-22"
    $cell.Style = "Normal"
}
